$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(39400, 2007, 3.145939949069287, 2008, 2.51031180018495),
    @(39583, 2008, 3.348624516477416, 2009, 6.253707197847636),
    @(39765, 2008, 1.769627576887389, 2009, -1.56363396419209),
    @(39948, 2009, -7.921319741078669, 2010, -14.35806537048445),
    @(40130, 2009, -4.774178217057779, 2010, 1.68591416918662),
    @(40310, 2010, 0.2827737265587604, 2011, 0.646969966709543),
    @(40494, 2010, 1.97975191822708, 2011, 3.609042024648068),
    @(40676, 2011, 4.453626355171969, 2012, 6.136355062499943),
    @(40862, 2011, 3.452886745653183, 2012, 2.11424984565185),
    @(41044, 2012, 1.477633171193093, 2013, 2.015050062499957),
    @(41228, 2012, 1.239479831392853, 2013, 0.9276342348636168),
    @(41409, 2013, -0.2005250704869344, 2014, 0.2533625145866525),
    @(41592, 2013, 0.2379616621361214, 2014, 1.3307042289459),
    @(41774, 2014, 2.234527904461125, 2015, 3.31205545981732),
    @(41957, 2014, 1.51977456621637, 2015, 0.3626364251072101),
    @(42137, 2015, 1.31172787261824, 2016, 1.093673275363716),
    @(42321, 2015, 1.470039379455756, 2016, 1.339087911421144),
    @(42503, 2016, 1.784808447869168, 2017, 2.687934870329323),
    @(42689, 2016, 1.638797242243251, 2017, 1.006353890555189),
    @(42867, 2017, 1.874682902292824, 2018, 2.445693358388845),
    @(43053, 2017, 2.161565493242668, 2018, 3.257358596620663),
    @(43145, 2018, 2.634902838428421, 2019, 2.450082126686781),
    @(43235, 2018, 2.159589514946769, 2019, 1.194430031759008),
    @(43326, 2018, 2.33065612535226, 2019, 1.797074113438568),
    @(43418, 2018, 2.214251681313772, 2019, -0.3013396321239648),
    @(43510, 2019, 0.2120367165967529, 2020, 0.06930121167902925),
    @(43600, 2019, 0.801449343663907, 2020, 1.609625625600009),
    @(43691, 2019, 0.5878492443567307, 2020, -0.1841852167622338),
    @(43783, 2019, 0.6066442151010376, 2020, 0.2691345740889695),
    @(43875, 2020, 0.2143710142614941, 2021, 0.112014348673295),
    @(43966, 2020, -3.107661574595755, 2021, -8.513835774399992),
    @(44068, 2020, -4.65090747647452, 2021, 16.69343145309163),
    @(44159, 2020, -4.207901339433196, 2021, -0.6332519459683494),
    @(44251, 2021, -0.3135161753049287, 2022, -2.807860884474733),
    @(44341, 2021, 0.6949587062036411, 2022, 4.530463903052695),
    @(44432, 2021, 0.7583924418458787, 2022, 5.353281688849409),
    @(44525, 2021, 1.099928004397532, 2022, 0.4297312830033428),
    @(44617, 2022, 2.258175284742725, 2023, 1.874171968961869),
    @(44706, 2022, 2.016988966764255, 2023, 1.492291178243965),
    @(44798, 2022, 2.236860175919531, 2023, 0.3985453946335094),
    @(44890, 2022, 2.310042359896225, 2023, -0.177017417229286),
    @(44981, 2023, -0.1441795889054687, 2024, -0.8289226263154981),
    @(45071, 2023, -0.1519071329076249, 2024, 0.3688835244738842),
    @(45163, 2023, -0.09588622947416248, 2024, -0.5974824858474559),
    @(45254, 2023, 0.0464415346324687, 2024, 0.08221134935635366),
    @(45345, 2024, -0.377538644385722, 2025, -0.3296319120892899),
    @(45436, 2024, -0.01173401322185352, 2025, 1.025808301409614),
    @(45534, 2024, -0.2385784141923808, 2025, -0.1293951476914468),
    @(45618, 2024, -0.3101476031197148, 2025, 0.5295174046934692),
    @(45713, 2025, 0.1482615292711431, 2026, 0.9611939328737495),
    @(45800, 2025, 0.1328390304517146, 2026, 0.2797831653477356),
    @(45891, 2025, -0.0960403240804597, 2026, 0.2498944033941575)
)

# Extend the table by one row (row 53), copying formatting from row 52
# so the new row's date cell (column A) keeps the existing date style.
$ws.Range("A52:E52").Copy()
$ws.Range("A53:E53").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value2 = $row[$c]
    }
}
